# Fix: Correct default plugin namespace; split namespace and keys with ':'
#
# Renames the three "excel:" metadata keys on the DataModels sheet to use
# the "excel.simple:" namespace, and restores the view state (active
# sheet/selection/column widths) to what a real save of the workbook after
# this edit looked like.

$wb = $excel.ActiveWorkbook
$dataModels = $wb.Worksheets.Item("DataModels")

# Rename the metadata keys to use the "excel.simple" namespace instead of
# "excel". These three header cells are the only cell content that
# actually changes.
$dataModels.Range("G1").Value = "uk.ac.ox.softeng.maurodatamapper.plugins.excel.simple:reviewed"
$dataModels.Range("H1").Value = "uk.ac.ox.softeng.maurodatamapper.plugins.excel.simple:approved"
$dataModels.Range("I1").Value = "uk.ac.ox.softeng.maurodatamapper.plugins.excel.simple:distributed"

# Widen the now-longer metadata-key columns to fit their new text.
$dataModels.Columns("G:H").ColumnWidth = 53.166666666666664
$dataModels.Columns("I:I").ColumnWidth = 54.45

# The DataModels sheet becomes the active tab/sheet, with a new scroll
# position and selected cell.
$dataModels.Activate()
$dataModels.Range("H3").Select()
